$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, [string]$val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '27.917.55'
$ws.Range("E2").Value = '  -2.55%  '
$ws.Range("D3").Value = '1.793.56'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("E4").Value = '  +0.13%  '
Set-TextValue "D5" '316.72'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("E6").Value = '  +0.12%  '
Set-TextValue "D7" '0.5304'
$ws.Range("E7").Value = '  -2.93%  '
Set-TextValue "D8" '0.3933'
$ws.Range("E8").Value = '  +3.23%  '
Set-TextValue "D9" '0.07452'
$ws.Range("E9").Value = '  -0.91%  '
Set-TextValue "D10" '41.50'
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("E11").Value = '  -2.84%  '
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D13" '7.483'
$ws.Range("E13").Value = '  +1.04%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D14" '6.171'
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("E15").Value = '  -1.97%  '
$ws.Range("D16").Value = '1.795.09'
$ws.Range("E16").Value = '  +0.23%  '
Set-TextValue "D17" '88.25'
$ws.Range("E17").Value = '  -2.28%  '
Set-TextValue "D18" '0.00001059'
$ws.Range("E18").Value = '  -0.85%  '
Set-TextValue "D19" '0.06586'
$ws.Range("E19").Value = '  +1.58%  '
$ws.Range("E20").Value = '  +0.07%  '
Set-TextValue "D21" '17.19'
$ws.Range("E21").Value = '  -0.93%  '
Set-TextValue "D22" '5.943'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '27.949.64'
$ws.Range("E23").Value = '  -2.36%  '
$ws.Range("E24").Value = '  -0.48%  '
Set-TextValue "D25" '2.094'
$ws.Range("E25").Value = '  +0.05%  '
Set-TextValue "D26" '156.89'
$ws.Range("E26").Value = '  -2.37%  '
$ws.Range("E27").Value = '  -1.57%  '
$ws.Range("D28").Value = '2.002.28'
$ws.Range("E28").Value = '  +0.05%  '
Set-TextValue "D29" '2.288'
$ws.Range("E29").Value = '  -3.30%  '
Set-TextValue "D30" '121.71'
$ws.Range("E30").Value = '  -1.32%  '
Set-TextValue "D31" '0.1084'
$ws.Range("E31").Value = '  +2.57%  '
Set-TextValue "D32" '1.094'
$ws.Range("E32").Value = '  -2.25%  '
Set-TextValue "D33" '3.677'
$ws.Range("E33").Value = '  -0.16%  '
Set-TextValue "D34" '5.493'
$ws.Range("E34").Value = '  -2.71%  '
Set-TextValue "D35" '0.07065'
$ws.Range("E35").Value = '  +6.14%  '
Set-TextValue "D36" '0.2206'
$ws.Range("E36").Value = '  -2.40%  '
Set-TextValue "D37" '5.099'
$ws.Range("E37").Value = '  +1.28%  '
Set-TextValue "D38" '0.02270'
$ws.Range("E38").Value = '  -1.41%  '
Set-TextValue "D39" '8.369'
$ws.Range("E39").Value = '  -4.66%  '
Set-TextValue "D40" '11.22'
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D41" '1.182'
$ws.Range("E41").Value = '  -1.28%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D42" '0.6110'
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("E43").Value = '  -1.35%  '
Set-TextValue "D44" '13.27'
$ws.Range("E44").Value = '  +0.14%  '
Set-TextValue "D45" '3.680'
$ws.Range("E45").Value = '  -0.45%  '
Set-TextValue "D46" '0.5701'
$ws.Range("E46").Value = '  -2.80%  '
Set-TextValue "D47" '125.28'
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("E48").Value = '  +1.73%  '
Set-TextValue "D49" '1.914'
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("E50").Value = '  -1.24%  '
Set-TextValue "D51" '71.17'
$ws.Range("E51").Value = '  -1.58%  '
